# rerun dist commute with harmonised education
# Applies the revised regression output (All_model_short / All_model_short (2))
# - corrects the UrbBuildDensity_res coefficient (B6) on both sheets
# - re-flags which p-values are shown in scientific notation on sheet 1
# - tidies formatting / row layout on sheet 2

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # All_model_short
$ws2 = $wb.Worksheets.Item(2)   # All_model_short (2)

# ---------------------------------------------------------------------------
# Sheet 1: All_model_short
# ---------------------------------------------------------------------------

# Corrected coefficient value for UrbBuildDensity_res
$ws1.Range("B6").Value = 67.0576977632157

# p-value column formatting: cells with a "significant" style keep the
# scientific-notation number format, the rest fall back to plain General
$ws1.Range("C2").NumberFormat  = "0.00E+00"
$ws1.Range("C3").NumberFormat  = "0.00E+00"
$ws1.Range("C4").NumberFormat  = "0.00E+00"
$ws1.Range("C5").Style         = "Normal"
$ws1.Range("C6").NumberFormat  = "0.00E+00"
$ws1.Range("C7").NumberFormat  = "0.00E+00"
$ws1.Range("C8").Style         = "Normal"
$ws1.Range("C9").NumberFormat  = "0.00E+00"
$ws1.Range("C10").Style        = "Normal"
$ws1.Range("C11").NumberFormat = "0.00E+00"
$ws1.Range("C12").Style        = "Normal"
$ws1.Range("C13").Style        = "Normal"
$ws1.Range("C14").Style        = "Normal"
$ws1.Range("C15").NumberFormat = "0.00E+00"

# Selection on sheet 1 moves to the coefficient column
$ws1.Activate()
$ws1.Range("B2:B15").Select()

# ---------------------------------------------------------------------------
# Sheet 2: All_model_short (2)
# ---------------------------------------------------------------------------

# Same corrected coefficient value
$ws2.Range("B6").Value = 67.0576977632157

# Coefficients now displayed with one decimal place instead of as integers
$ws2.Range("B5").NumberFormat  = "0.0"
$ws2.Range("B6").NumberFormat  = "0.0"
$ws2.Range("B7").NumberFormat  = "0.0"
$ws2.Range("B8").NumberFormat  = "0.0"
$ws2.Range("B11").NumberFormat = "0.0"
$ws2.Range("B12").NumberFormat = "0.0"
$ws2.Range("B13").NumberFormat = "0.0"
$ws2.Range("B14").NumberFormat = "0.0"
$ws2.Range("B15").NumberFormat = "0.0"

# Row 4 label loses its special formatting (back to default)
$ws2.Range("A4").Style = "Normal"

# Rows 5-17 labels are harmonised onto a single "vertical-center" style
$ws2.Range("A5:A17").VerticalAlignment = -4108   # xlVAlignCenter

# Remove the two trailing blank rows (18 and 19)
$ws2.Rows("18:19").Delete()

# Selection on sheet 2 moves, sheet 2 remains the active tab
$ws2.Activate()
$ws2.Range("B6").Select()
